$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")
$ws.Rows("2:3").Insert()
$ws.Range("A2:F3").ClearFormats()

$ws.Range("A2").Value = "'3975"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = "'6"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'30.94%"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").Value = "'3977"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = "NO"

$ws.Range("A2:F3").Style = "Normal"
